$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Shai Gilgeous-Alexander",
    "Jaylen Clark",
    "Cason Wallace",
    "Jamal Murray",
    "Royce O'Neale",
    "Bam Adebayo",
    "Kris Dunn",
    "CJ McCollum",
    "Tobias Harris",
    "Kyle Kuzma",
    "John Collins",
    "Jordan Clarkson",
    "Kyshawn George",
    "Kyrie Irving",
    "Zach LaVine",
    "Rui Hachimura",
    "Lauri Markkanen",
    "Jordan Poole"
)

$positions = @(
    "PG,SG",
    "SG",
    "PG,SG",
    "PG,SG",
    "SF,PF",
    "PF,C",
    "PG,SG",
    "PG,SG",
    "SF,PF",
    "SF,PF",
    "PF,C",
    "SG,SF",
    "SG,SF",
    "PG,SG",
    "SG,SF",
    "SF,PF",
    "SF,PF",
    "PG,SG"
)

$teams = @(
    "Oklahoma City Thunder",
    "Minnesota Timberwolves",
    "Oklahoma City Thunder",
    "Denver Nuggets",
    "Phoenix Suns",
    "Miami Heat",
    "LA Clippers",
    "New Orleans Pelicans",
    "Detroit Pistons",
    "Milwaukee Bucks",
    "Utah Jazz",
    "Utah Jazz",
    "Washington Wizards",
    "Dallas Mavericks",
    "Sacramento Kings",
    "Los Angeles Lakers",
    "Utah Jazz",
    "Washington Wizards"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $positions[$i]
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
